# Generate Report for Handback
# The file "a788d847-8a9b-4446-992f-1622740ddb20.md" has now been handed back
# (in sync with en-US) instead of merely "Ready for handoff". Update the
# Overview sheet as well as the per-language (zh-cn / de-de) detail sheets to
# reflect the new status, the refreshed "Latest Handback DateTime" timestamps,
# and clear the stale "Error Detail" message that no longer applies.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: update the zh-cn / de-de status columns for the row that
# corresponds to a788d847-8a9b-4446-992f-1622740ddb20.md (row 3).
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $handedBack
$overview.Range("F3").Value = $handedBack

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 is the a788d847... file. Status moves to "handed back",
# the handback datetime is refreshed, and the error detail is cleared out.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $handedBack
$zhcn.Range("K3").Value = "2016-08-20 18:57:29"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 13

# ---------------------------------------------------------------------------
# de-de sheet: same update, with its own refreshed handback datetime.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $handedBack
$dede.Range("K3").Value = "2016-08-20 18:57:35"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 13
